# Daily attendance processing - 2026-01-27 21:33:43
# Normalize the "Recorded By" (column G) ordering: swap
# "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

$oldText = "dnasr281@gmail.com, System"
$newText = "System, dnasr281@gmail.com"

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Value2 -eq $oldText) {
        $cell.Value = $newText
    }
}
